$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of rows 2 and 3 entirely (spans A:AD), leaving them blank
# without shifting row 4 up.
$ws.Range("A2:AD3").ClearContents()

# Update row 4 with the new values
$ws.Range("A4").Value = 0.3296973392106938
$ws.Range("B4").Value = 0.417805974292161
$ws.Range("C4").Value = 0.2413099029310042
$ws.Range("D4").Value = 0.3735848048328873
$ws.Range("E4").Value = 0.2840026953544471
$ws.Range("F4").Value = 7.002077942463496
$ws.Range("G4").Value = 8.695149276628957
$ws.Range("H4").Value = 5.310634041337306
$ws.Range("I4").Value = 7.843189284973572
$ws.Range("J4").Value = 6.125634627405966
$ws.Range("K4").Value = 2.308718684912954
$ws.Range("L4").Value = 3.198881118013062
$ws.Range("M4").Value = 1.530571468379241
$ws.Range("N4").Value = 2.723243821056452
$ws.Range("O4").Value = 1.874636325871613
$ws.Range("P4").Value = 67.25612
$ws.Range("Q4").Value = 114.6473069500463
$ws.Range("R4").Value = 43.08260692625805
$ws.Range("S4").Value = 81.81154969567304
$ws.Range("T4").Value = 51.47837976500801
$ws.Range("U4").Value = 0.1997709064090712
$ws.Range("V4").Value = 0.3236004793987665
$ws.Range("W4").Value = 0.06896070333280385
$ws.Range("X4").Value = 0.2665426857909344
$ws.Range("Y4").Value = 0.1280229037497588
$ws.Range("Z4").Value = 0.8378709402063405
$ws.Range("AA4").Value = 0.96809666103918
$ws.Range("AB4").Value = 0.5926709557689795
$ws.Range("AC4").Value = 0.9080364404473633
$ws.Range("AD4").Value = 0.7360443419493188
